# Auto-generated Excel COM-interop script to apply cryptos list update
# Commit: Updated cryptos list on Sun Dec 17 11:25:46 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.042.84"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "'2.217.90"
$ws.Range("E3").Value = "  -1.27%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").Value = "'242.09"
$ws.Range("E5").Value = "  -1.70%  "

$ws.Range("D6").Value = "'0.628"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").Value = "'73.16"
$ws.Range("E7").Value = "  -2.09%  "

$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("D10").Value = "'42.62"
$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").Value = "'0.0959"
$ws.Range("E11").Value = "  +1.24%  "

$ws.Range("D12").Value = "'7.07"
$ws.Range("E12").Value = "  -1.39%  "

$ws.Range("E13").Value = "  +0.31%  "

$ws.Range("D14").Value = "'2.551.23"
$ws.Range("E14").Value = "  -1.24%  "

$ws.Range("D15").Value = "'14.31"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("D16").Value = "'0.839"
$ws.Range("E16").Value = "  -1.77%  "

$ws.Range("D17").Value = "'2.206.67"
$ws.Range("E17").Value = "  -1.20%  "

$ws.Range("D18").Value = "'41.913.26"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("E19").Value = "  +6.91%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "'72.77"

$ws.Range("D22").Value = "'10.79"
$ws.Range("E22").Value = "  +18.61%  "

$ws.Range("D23").Value = "'230.60"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("D24").Value = "'2.08"
$ws.Range("E24").Value = "  -5.89%  "

$ws.Range("D25").Value = "'11.86"
$ws.Range("E25").Value = "  +3.78%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").Value = "'3.67"
$ws.Range("E27").Value = "  +1.28%  "

$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  -1.36%  "

$ws.Range("E29").Value = "  -3.02%  "

$ws.Range("D30").Value = "'168.13"
$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("D31").Value = "'20.51"
$ws.Range("E31").Value = "  -0.84%  "

$ws.Range("D32").Value = "'5.68"
$ws.Range("E32").Value = "  +8.98%  "

$ws.Range("D33").Value = "'0.0797"
$ws.Range("E33").Value = "  -3.08%  "

$ws.Range("D34").Value = "'30.08"
$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("D36").Value = "'0.110"
$ws.Range("E36").Value = "  -10.30%  "

$ws.Range("D37").Value = "'4.28"
$ws.Range("E37").Value = "  -3.84%  "

$ws.Range("E38").Value = "  -4.28%  "

$ws.Range("D39").Value = "'13.92"
$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("D40").Value = "'65.60"
$ws.Range("E40").Value = "  +4.96%  "

$ws.Range("D41").Value = "'2.13"
$ws.Range("E41").Value = "  -2.50%  "

$ws.Range("D42").Value = "'5.65"
$ws.Range("E42").Value = "  -2.48%  "

$ws.Range("D43").Value = "'0.199"
$ws.Range("E43").Value = "  -2.62%  "

$ws.Range("D44").Value = "'8.81"
$ws.Range("E44").Value = "  +0.95%  "

$ws.Range("D45").Value = "'105.70"
$ws.Range("E45").Value = "  -1.95%  "

$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("E47").Value = "  +6.10%  "

$ws.Range("D48").Value = "'1.12"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("E50").Value = "  -0.22%  "

$ws.Range("D51").Value = "'2.422.01"
$ws.Range("E51").Value = "  -1.56%  "
